$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert two new columns at D:E ---
# This shifts existing columns D:K to F:M, preserving their values/styles.
$ws.Columns("D:E").Insert()

# --- Step 2: Copy formatting (number format/style) from column F (the old column D)
#     into the two newly inserted blank columns D:E, so they render the same way
#     (date format row 7/38/80, thousands format elsewhere). ---
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 3: Populate the new D:E columns with the new quarterly figures ---
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 28600
$ws.Range("E8").Value2 = 17900
$ws.Range("D9").Value2 = 22600
$ws.Range("E9").Value2 = 16300
$ws.Range("D10").Value2 = 6000
$ws.Range("E10").Value2 = 1600
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 25300
$ws.Range("E17").Value2 = 19000
$ws.Range("D18").Value2 = 3300
$ws.Range("E18").Value2 = -1100
$ws.Range("D20").Value2 = 800
$ws.Range("E20").Value2 = 300
$ws.Range("D21").Value2 = 6600
$ws.Range("E21").Value2 = 1400
$ws.Range("D22").Value2 = 2400
$ws.Range("E22").Value2 = 2200
$ws.Range("D23").Value2 = 1800
$ws.Range("E23").Value2 = -2900
$ws.Range("D24").Value2 = 500
$ws.Range("E24").Value2 = -500
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 1300
$ws.Range("E26").Value2 = -2400
$ws.Range("D27").Value2 = 1300
$ws.Range("E27").Value2 = -2400
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = -200
$ws.Range("E29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -800
$ws.Range("E32").Value2 = -300
$ws.Range("D33").Value2 = 1100
$ws.Range("E33").Value2 = -2400
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 1100
$ws.Range("E35").Value2 = -2400
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 19000
$ws.Range("E41").Value2 = 21200
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 0
$ws.Range("E43").Value2 = 0
$ws.Range("D44").Value2 = 177100
$ws.Range("E44").Value2 = 179600
$ws.Range("D45").Value2 = 0
$ws.Range("E45").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("E46").Value2 = 0
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 253100
$ws.Range("E48").Value2 = 234800
$ws.Range("D49").Value2 = 0
$ws.Range("E49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 31700
$ws.Range("E52").Value2 = 38500
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 496500
$ws.Range("E54").Value2 = 488100
$ws.Range("D57").Value2 = 20600
$ws.Range("E57").Value2 = 21000
$ws.Range("D58").Value2 = 0
$ws.Range("E58").Value2 = 0
$ws.Range("D59").Value2 = 11900
$ws.Range("E59").Value2 = 10400
$ws.Range("D60").Value2 = 0
$ws.Range("E60").Value2 = 0
$ws.Range("D61").Value2 = 295500
$ws.Range("E61").Value2 = 293700
$ws.Range("D62").Value2 = 0
$ws.Range("E62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 372500
$ws.Range("E66").Value2 = 365400
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = -41100
$ws.Range("E72").Value2 = -42200
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 124000
$ws.Range("E76").Value2 = 122600
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 1100
$ws.Range("E81").Value2 = -2400
$ws.Range("D83").Value2 = 2400
$ws.Range("E83").Value2 = 2200
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = -5300
$ws.Range("E89").Value2 = -17900
$ws.Range("D91").Value2 = -8500
$ws.Range("E91").Value2 = -10500
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -9500
$ws.Range("E94").Value2 = -10600
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = 6600
$ws.Range("E100").Value2 = 38200
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("D102").Value2 = -8200
$ws.Range("E102").Value2 = 9800

# --- Step 4: Apply individual restated-value corrections to the shifted F:M columns ---
$ws.Range("I8").Value2 = 59300
$ws.Range("I9").Value2 = 51300
$ws.Range("I10").Value2 = 8000
$ws.Range("I17").Value2 = 37400
$ws.Range("I18").Value2 = 21900
$ws.Range("I20").Value2 = 200
$ws.Range("I21").Value2 = 28000
$ws.Range("I22").Value2 = 5100
$ws.Range("I23").Value2 = 17000
$ws.Range("H24").Value2 = 100
$ws.Range("I24").Value2 = 6200
$ws.Range("H26").Value2 = 700
$ws.Range("I26").Value2 = 10800
$ws.Range("H27").Value2 = 700
$ws.Range("I27").Value2 = 10700
$ws.Range("H29").Value2 = -7600
$ws.Range("I32").Value2 = -200
$ws.Range("I33").Value2 = 10700
$ws.Range("I35").Value2 = 10700
$ws.Range("I81").Value2 = 10700
$ws.Range("I83").Value2 = 5900
$ws.Range("H89").Value2 = 25700
$ws.Range("I89").Value2 = -12900
$ws.Range("I91").Value2 = -14400
$ws.Range("J91").Value2 = -2800
$ws.Range("H94").Value2 = -20200
$ws.Range("I94").Value2 = 101200
$ws.Range("I96").Value2 = -8100
$ws.Range("I100").Value2 = -73900
$ws.Range("H102").Value2 = -500
$ws.Range("I102").Value2 = 14400
